# Capitalize the hex letters (A-F) in the byte-code strings found in
# columns G (doip) and H (uds) of the active worksheet.
# Values look like "0x02:0xfd:0x00:..." and must become "0x02:0xFD:0x00:...".
# The leading "0x" marker itself must stay lowercase; only the two hex
# digits following it are upper-cased. Non hex values (e.g. "N/A") are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Convert-HexUpper($s) {
    if ($null -eq $s) { return $s }
    $parts = $s.Split(":")
    $result = @()
    foreach ($p in $parts) {
        if ($p.Length -ge 2 -and $p.Substring(0, 2) -eq "0x") {
            $result += "0x" + $p.Substring(2).ToUpper()
        } else {
            $result += $p
        }
    }
    return [string]::Join(":", $result)
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    foreach ($col in @("G", "H")) {
        $cell = $ws.Range($col + $r)
        $val = $cell.Value2
        if ($val -is [string]) {
            $newVal = Convert-HexUpper $val
            # Always (re)assign; the conversion is idempotent for values
            # that contain no "0x.." hex tokens (e.g. "N/A"), and the
            # string comparison operators in this engine are
            # case-insensitive, so we cannot reliably skip unchanged cells.
            $cell.Value2 = $newVal
        }
    }
}
